$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are stored as text so numeric-looking strings
# (e.g. "1.001", "28.052.36") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.052.36"
$ws.Range("E2").Value = "  -3.44%  "

$ws.Range("D3").Value = "1.748.31"
$ws.Range("E3").Value = "  -3.97%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Value = "226.33"
$ws.Range("E5").Value = "  -3.21%  "

$ws.Range("D6").Value = "0.5806"
$ws.Range("E6").Value = "  -2.89%  "

$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").Value = "0.2717"
$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("D9").Value = "23.28"
$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("D10").Value = "0.06611"
$ws.Range("E10").Value = "  -4.40%  "

$ws.Range("D11").Value = "0.07501"
$ws.Range("E11").Value = "  -1.10%  "

$ws.Range("D12").Value = "1.748.00"
$ws.Range("E12").Value = "  -4.28%  "

$ws.Range("D13").Value = "4.714"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").Value = "0.6062"
$ws.Range("E14").Value = "  -2.52%  "

$ws.Range("D15").Value = "1.980.81"
$ws.Range("E15").Value = "  -4.24%  "

$ws.Range("D16").Value = "74.06"
$ws.Range("E16").Value = "  -3.82%  "

$ws.Range("D17").Value = "0.000008650"
$ws.Range("E17").Value = "  -10.33%  "

$ws.Range("D18").Value = "28.015.10"
$ws.Range("E18").Value = "  -2.33%  "

$ws.Range("D19").Value = "5.322"
$ws.Range("E19").Value = "  -4.13%  "

$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").Value = "205.79"
$ws.Range("E21").Value = "  -4.30%  "

$ws.Range("D22").Value = "11.30"
$ws.Range("E22").Value = "  -1.54%  "

$ws.Range("D23").Value = "6.648"
$ws.Range("E23").Value = "  -2.75%  "

$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").Value = "150.50"
$ws.Range("E25").Value = "  -3.68%  "

$ws.Range("D26").Value = "8.023"
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("D27").Value = "0.1238"
$ws.Range("E27").Value = "  -3.02%  "

$ws.Range("D28").Value = "16.12"
$ws.Range("E28").Value = "  -1.44%  "

$ws.Range("D29").Value = "1.394"
$ws.Range("E29").Value = "  -1.67%  "

$ws.Range("D30").Value = "0.06107"
$ws.Range("E30").Value = "  -4.17%  "

$ws.Range("D31").Value = "1.386"
$ws.Range("E31").Value = "  -3.42%  "

$ws.Range("D32").Value = "3.743"
$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("D33").Value = "3.722"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("D34").Value = "1.686"
$ws.Range("E34").Value = "  -1.76%  "

$ws.Range("D35").Value = "1.038"
$ws.Range("E35").Value = "  -4.20%  "

$ws.Range("D36").Value = "0.6362"
$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("D37").Value = "2.429"
$ws.Range("E37").Value = "  -4.17%  "

$ws.Range("D38").Value = "2.653"
$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("D39").Value = "0.01671"
$ws.Range("E39").Value = "  -4.31%  "

$ws.Range("D40").Value = "6.272"
$ws.Range("E40").Value = "  -3.63%  "

$ws.Range("D41").Value = "1.129.57"
$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("D42").Value = "0.8641"
$ws.Range("E42").Value = "  -1.72%  "

$ws.Range("D43").Value = "1.005"
$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("D44").Value = "99.73"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("D45").Value = "1.890.24"
$ws.Range("E45").Value = "  -4.25%  "

$ws.Range("D46").Value = "59.17"
$ws.Range("E46").Value = "  -3.52%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "1.576"
$ws.Range("E47").Value = "  -1.28%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.00000000107"
$ws.Range("E48").Value = "  -4.31%  "

$ws.Range("D49").Value = "8.245"
$ws.Range("E49").Value = "  -2.24%  "

$ws.Range("D50").Value = "0.05398"
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("D51").Value = "0.4437"
$ws.Range("E51").Value = "  -1.96%  "

# Restore the default "Normal" style on column D so no residual text
# number-format is left applied to the cells (matches original styling).
$ws.Range("D2:D51").Style = "Normal"